$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "iNeuron.ai"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "Caster"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "345"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "456"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "24"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "56"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "234"
